$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting all columns (B:F) one position to the left (-> A:E)
$ws.Columns.Item(1).Delete()
